# T7, T9 and ongoing JS work — add new glyph row (g34 / pro abbreviation)
# to the "Glyphs" sheet, right after the existing g33 / Vertical m row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Glyphs")
$ws.Activate()

$ws.Range("A35").Value = "g34"
$ws.Range("B35").Value = "pro abbreviation"

# Mirror the author's saved view-state: selection moved on to the newly
# entered cell.
$ws.Range("B35").Select() | Out-Null
